$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 800, shifting existing rows 800-841 down to 801-842
$ws.Rows.Item(800).Insert()

# Populate the newly inserted row 800 with the new data.
# Force column A to be treated as plain text (not auto-parsed into a date
# serial number) and then reset the cell style back to the sheet's default
# so no stray number-format style is introduced.
$ws.Cells.Item(800, 1).NumberFormat = "@"
$ws.Cells.Item(800, 1).Value2 = "2026/02/15"
$ws.Cells.Item(800, 1).Style = "Normal"

$ws.Cells.Item(800, 2).Value = "日"
$ws.Cells.Item(800, 3).Value = 19
$ws.Cells.Item(800, 4).Value = 31
